$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 20 updates: Oct 1 work session ran later than first logged, until 23:45 ---
# D20: end time changes from 21:45 to 23:45 -> reuse the standard time-format
# style from C20 (the cell had picked up a one-off bold/italic AM-PM style;
# normalize it back to the shared style used by the rest of column D)
$ws.Range("C20").Copy()
$ws.Range("D20").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D20").Value = $excel.Evaluate("TIME(23,45,0)")

# G20: add progress note, reusing the note style already used in G18/G19
$ws.Range("G18").Copy()
$ws.Range("G20").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G20").Value = "Section 5, finished upto lesson# 46"

# --- Row 21: new work session logged for Oct 2 ---
$ws.Range("B21").Value = $excel.Evaluate("DATE(2022,10,2)")
$ws.Range("C21").Value = $excel.Evaluate("TIME(20,30,0)")
$ws.Range("D21").Value = $excel.Evaluate("TIME(21,30,0)")
$ws.Range("E21").Formula = "=D21-C21"

# G21: progress note for this session, same style as the other note cells
$ws.Range("G18").Copy()
$ws.Range("G21").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G21").Value = "Section 5, finished upto lesson# 48"

# Leave the cursor on the newly-entered note, matching the saved view state
$ws.Range("G21").Select()

$wb.Save()
